$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117 (shifts existing rows 117-206 down to 118-207)
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record's data
$ws.Cells.Item(117, 1).Value = 11
$ws.Cells.Item(117, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(117, 3).Value = "Bíobío"
$ws.Cells.Item(117, 4).Value = 45090
$ws.Cells.Item(117, 5).Value = 8
$ws.Cells.Item(117, 6).Value = 100112021
$ws.Cells.Item(117, 7).Value = "Ají"
$ws.Cells.Item(117, 8).Value = "Inferno"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 50
$ws.Cells.Item(117, 11).Value = 13000
$ws.Cells.Item(117, 12).Value = 14000
$ws.Cells.Item(117, 13).Value = 13400
$ws.Cells.Item(117, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(117, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(117, 16).Value = 1340
$ws.Cells.Item(117, 17).Value = 10
$ws.Cells.Item(117, 18).Value = "Hortaliza"
